$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; existing rows 18-21 shift down to 19-22.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 44984
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100114007
$ws.Cells.Item(18, 7).Value = "Jengibre"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 16000
$ws.Cells.Item(18, 12).Value = 17000
$ws.Cells.Item(18, 13).Value = 16500
$ws.Cells.Item(18, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 1269
$ws.Cells.Item(18, 17).Value = 13
$ws.Cells.Item(18, 18).Value = "Hortaliza"
